$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1856987.6
$ws.Range("J17").Value = 1856987.6
$ws.Range("L17").Value = 5570962.800000001
$ws.Range("N17").Value = -5571298.800000001
$ws.Range("H38").Value = 636.8
$ws.Range("I38").Value = 53
$ws.Range("J38").Value = 2972
$ws.Range("K38").Value = 159
$ws.Range("L38").Value = 8916
$ws.Range("M38").Value = 213
$ws.Range("N38").Value = -9660
$ws.Range("H58").Value = 652.6
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300
$ws.Range("H62").Value = 4294.6875
$ws.Range("I62").Value = 3271.5
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 3271.5
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -2647.5
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 4294.6875
$ws.Range("I65").Value = 3271.5
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 16357.5
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -13237.5
$ws.Range("N65").Value = -36240
$ws.Range("H87").Value = 15306.411
$ws.Range("J87").Value = 15306.411
$ws.Range("L87").Value = 15306.411
$ws.Range("N87").Value = -17802.411
$ws.Range("H90").Value = 15306.411
$ws.Range("J90").Value = 15306.411
$ws.Range("L90").Value = 45919.233
$ws.Range("N90").Value = -58399.233
$ws.Range("H129").Value = 618.38464
$ws.Range("I129").Value = 503.9
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 1511.7
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 3488.3
$ws.Range("N129").Value = -13000
$ws.Range("H138").Value = 6761.5073
$ws.Range("I138").Value = 4699
$ws.Range("J138").Value = 6964.377
$ws.Range("K138").Value = 14097
$ws.Range("L138").Value = 20893.131
$ws.Range("M138").Value = -8957
$ws.Range("N138").Value = -31173.131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19072.215
$ws.Range("I32").Value = 12876.5
$ws.Range("J32").Value = 27333.166
$ws.Range("K32").Value = 12876.5
$ws.Range("L32").Value = 27333.166
$ws.Range("M32").Value = -12589.5
$ws.Range("N32").Value = -27907.166
$ws.Range("H61").Value = 1919.5
$ws.Range("I61").Value = 1642.1428
$ws.Range("J61").Value = 2566.6667
$ws.Range("K61").Value = 1642.1428
$ws.Range("L61").Value = 2566.6667
$ws.Range("M61").Value = -1430.1428
$ws.Range("N61").Value = -2990.6667
$ws.Range("H102").Value = 1691.8182
$ws.Range("I102").Value = 1374.1428
$ws.Range("J102").Value = 2247.75
$ws.Range("K102").Value = 1374.1428
$ws.Range("L102").Value = 2247.75
$ws.Range("M102").Value = 247.8571999999999
$ws.Range("N102").Value = -5491.75
$ws.Range("H122").Value = 2354.1
$ws.Range("I122").Value = 1805.2858
$ws.Range("J122").Value = 3634.6667
$ws.Range("K122").Value = 5415.857400000001
$ws.Range("L122").Value = 10904.0001
$ws.Range("M122").Value = -2965.857400000001
$ws.Range("N122").Value = -15804.0001
$ws.Range("H132").Value = 1683.875
$ws.Range("I132").Value = 1290.7916
$ws.Range("J132").Value = 2863.125
$ws.Range("K132").Value = 3872.3748
$ws.Range("L132").Value = 8589.375
$ws.Range("M132").Value = -1342.3748
$ws.Range("N132").Value = -13649.375
$ws.Range("H136").Value = 1919.5
$ws.Range("I136").Value = 1642.1428
$ws.Range("J136").Value = 2566.6667
$ws.Range("K136").Value = 4926.428400000001
$ws.Range("L136").Value = 7700.000100000001
$ws.Range("M136").Value = -2376.428400000001
$ws.Range("N136").Value = -12800.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 20007500
$ws.Range("J19").Value = 20007500
$ws.Range("L19").Value = 20007500
$ws.Range("N19").Value = -20007846
$ws.Range("H134").Value = 2661.0908
$ws.Range("I134").Value = 2014.4
$ws.Range("J134").Value = 3200
$ws.Range("K134").Value = 6043.200000000001
$ws.Range("L134").Value = 9600
$ws.Range("M134").Value = -3508.200000000001
$ws.Range("N134").Value = -14670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1125.2106
$ws.Range("I16").Value = 1075.4615
$ws.Range("J16").Value = 1233
$ws.Range("K16").Value = 1075.4615
$ws.Range("L16").Value = 1233
$ws.Range("M16").Value = -788.4614999999999
$ws.Range("N16").Value = -1807
$ws.Range("H31").Value = 8746.375
$ws.Range("I31").Value = 3421.8572
$ws.Range("J31").Value = 16200.7
$ws.Range("K31").Value = 3421.8572
$ws.Range("L31").Value = 16200.7
$ws.Range("M31").Value = -3126.8572
$ws.Range("N31").Value = -16790.7
$ws.Range("H34").Value = 8746.375
$ws.Range("I34").Value = 3421.8572
$ws.Range("J34").Value = 16200.7
$ws.Range("K34").Value = 3421.8572
$ws.Range("L34").Value = 16200.7
$ws.Range("M34").Value = -3219.8572
$ws.Range("N34").Value = -16604.7
$ws.Range("H43").Value = 60000
$ws.Range("J43").Value = 60000
$ws.Range("L43").Value = 60000
$ws.Range("N43").Value = -60368
$ws.Range("H58").Value = 2778.75
$ws.Range("I58").Value = 1984
$ws.Range("J58").Value = 3573.5
$ws.Range("K58").Value = 1984
$ws.Range("L58").Value = 3573.5
$ws.Range("M58").Value = -1781
$ws.Range("N58").Value = -3979.5
$ws.Range("H101").Value = 60000
$ws.Range("J101").Value = 60000
$ws.Range("L101").Value = 60000
$ws.Range("N101").Value = -66490
$ws.Range("H113").Value = 1125.2106
$ws.Range("I113").Value = 1075.4615
$ws.Range("J113").Value = 1233
$ws.Range("K113").Value = 1075.4615
$ws.Range("L113").Value = 1233
$ws.Range("M113").Value = 1094.5385
$ws.Range("N113").Value = -5573
$ws.Range("H132").Value = 3022.647
$ws.Range("I132").Value = 1645.7778
$ws.Range("J132").Value = 4571.625
$ws.Range("K132").Value = 4937.3334
$ws.Range("L132").Value = 13714.875
$ws.Range("M132").Value = -2407.3334
$ws.Range("N132").Value = -18774.875
$ws.Range("H136").Value = 2778.75
$ws.Range("I136").Value = 1984
$ws.Range("J136").Value = 3573.5
$ws.Range("K136").Value = 5952
$ws.Range("L136").Value = 10720.5
$ws.Range("M136").Value = -3402
$ws.Range("N136").Value = -15820.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 8016.9395
$ws.Range("I137").Value = 2922.7144
$ws.Range("J137").Value = 9388.462
$ws.Range("K137").Value = 8768.143199999999
$ws.Range("L137").Value = 28165.386
$ws.Range("M137").Value = -3668.143199999999
$ws.Range("N137").Value = -38365.386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 904.56525
$ws.Range("I2").Value = 699.6667
$ws.Range("K2").Value = 699.6667
$ws.Range("M2").Value = -586.6667
$ws.Range("H70").Value = 6615.5713
$ws.Range("I70").Value = 4166.6665
$ws.Range("K70").Value = 4166.6665
$ws.Range("M70").Value = -3896.6665
$ws.Range("H73").Value = 6615.5713
$ws.Range("I73").Value = 4166.6665
$ws.Range("K73").Value = 4166.6665
$ws.Range("M73").Value = -3230.6665
$ws.Range("H132").Value = 2932.0303
$ws.Range("I132").Value = 2602.25
$ws.Range("J132").Value = 3439.3845
$ws.Range("K132").Value = 7806.75
$ws.Range("L132").Value = 10318.1535
$ws.Range("M132").Value = -5276.75
$ws.Range("N132").Value = -15378.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 801
$ws.Range("I93").Value = 801
$ws.Range("K93").Value = 801
$ws.Range("M93").Value = 447
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H122").Value = 27785278
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 9048.27
$ws.Range("I132").Value = 3268.1177
$ws.Range("J132").Value = 19966.334
$ws.Range("K132").Value = 9804.3531
$ws.Range("L132").Value = 59899.00199999999
$ws.Range("M132").Value = -7274.3531
$ws.Range("N132").Value = -64959.00199999999
$ws.Range("H136").Value = 47624960
$ws.Range("I136").Value = 7676.6
$ws.Range("J136").Value = 166668180
$ws.Range("K136").Value = 23029.8
$ws.Range("L136").Value = 500004540
$ws.Range("M136").Value = -20479.8
$ws.Range("N136").Value = -500009640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 766.6667
$ws.Range("I100").Value = 720
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1440
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -899
$ws.Range("N100").Value = -3082
$ws.Range("H126").Value = 86048.766
$ws.Range("I126").Value = 101121.27
$ws.Range("J126").Value = 3150
$ws.Range("K126").Value = 303363.81
$ws.Range("L126").Value = 9450
$ws.Range("M126").Value = -300893.81
$ws.Range("N126").Value = -14390
$ws.Range("H132").Value = 1969.5077
$ws.Range("I132").Value = 1411.9767
$ws.Range("K132").Value = 4235.9301
$ws.Range("M132").Value = -1705.9301
$ws.Range("H136").Value = 5747.9624
$ws.Range("I136").Value = 4234
$ws.Range("J136").Value = 6908.6665
$ws.Range("K136").Value = 12702
$ws.Range("L136").Value = 20725.9995
$ws.Range("M136").Value = -10152
$ws.Range("N136").Value = -25825.9995
